$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (D = date serial, M = Volumen, N = Precio minimo,
# O = Precio maximo, P = Precio promedio ponderado, Q = Unidad de
# comercializacion, S = Precio $/Kg). Columns A,B,C,E,F,G,H,I,J,K,L,R,T
# are unchanged.

$rowsData = @(
    @{ Row = 2; D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ Row = 3; D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 4; D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";           S = 786 },
    @{ Row = 5; D = 44491; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 14 kilos empedrada"; S = 643 },
    @{ Row = 6; D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";           S = 714 },
    @{ Row = 7; D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 8; D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 9; D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 }
)

foreach ($item in $rowsData) {
    $r = $item.Row
    $ws.Range("D$r").Value = $item.D
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
    $ws.Range("O$r").Value = $item.O
    $ws.Range("P$r").Value = $item.P
    $ws.Range("Q$r").Value = $item.Q
    $ws.Range("S$r").Value = $item.S
}
